$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 12 new columns starting at B (old B:J shift right to N:V)
$ws.Range("B1:M1").EntireColumn.Insert()

# Fill the newly inserted B1:M1 with the repeating pattern 1,2,3
$pattern = 1,2,3,1,2,3,1,2,3,1,2,3
for ($i = 0; $i -lt 12; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $pattern[$i]
}

# Match the width used by column A across the newly inserted columns
$ws.Range("B1:M1").EntireColumn.ColumnWidth = 14

# Update the view: scrolled to show column G first, last column selected
$excel.ActiveWindow.ScrollColumn = 7
[void]$ws.Range("V1:V1048576").Select()

Write-Output "done"
